# ---------------------------------------------------------------------------
# chore: update Sheets via scheduled runner
#
# Refreshes the market-board snapshot (currentAveragePrice / NQ / HQ) and the
# derived LevePrice / LeveProfit columns (H:N) for the affected rows in each
# per-job Leve-profit table (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# A few rows flip whether an HQ/NQ profit is even meaningful once the refreshed
# price makes LevePriceNQ/LevePriceHQ (K/L) go to - or come off of - zero, so the
# corresponding LeveProfitNQ/LeveProfitHQ cell (M/N) is created or cleared rather
# than just renumbered.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ==== ALC sheet ====
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Cells.Item(2, 8).Value = 600.8570999999999  # H2: 739.7273 -> 600.8570999999999
$ws.Cells.Item(2, 9).Value = 289.72726  # I2: 313.7 -> 289.72726
$ws.Cells.Item(2, 10).Value = 1741.6666  # J2: 5000 -> 1741.6666
$ws.Cells.Item(2, 11).Value = 289.72726  # K2: 313.7 -> 289.72726
$ws.Cells.Item(2, 12).Value = 1741.6666  # L2: 5000 -> 1741.6666
$ws.Cells.Item(2, 13).Value = -176.72726  # M2: -200.7 -> -176.72726
$ws.Cells.Item(2, 14).Value = -1967.6666  # N2: -5226 -> -1967.6666
# Row 15 (Leve Item ID 44146)
$ws.Cells.Item(15, 8).Value = 2378.5085  # H15: 2390.224 -> 2378.5085
$ws.Cells.Item(15, 9).Value = 2378.5085  # I15: 2390.224 -> 2378.5085
$ws.Cells.Item(15, 11).Value = 7135.5255  # K15: 7170.672 -> 7135.5255
$ws.Cells.Item(15, 13).Value = -6966.5255  # M15: -7001.672 -> -6966.5255
# Row 43 (Leve Item ID 5472)
$ws.Cells.Item(43, 8).Value = 6833  # H43: 9500 -> 6833
$ws.Cells.Item(43, 9).Value = 1499  # I43: 0 -> 1499
$ws.Cells.Item(43, 11).Value = 1499  # K43: 0 -> 1499
$ws.Cells.Item(43, 13).Value = -1430  # M43: None -> -1430
# Row 106 (Leve Item ID 19903)
$ws.Cells.Item(106, 8).Value = 8235444.5  # H106: 7720760 -> 8235444.5
$ws.Cells.Item(106, 9).Value = 12349678  # I106: 11227025 -> 12349678
$ws.Cells.Item(106, 11).Value = 12349678  # K106: 11227025 -> 12349678
$ws.Cells.Item(106, 13).Value = -12349047  # M106: -11226394 -> -12349047
# Row 132 (Leve Item ID 44049)
$ws.Cells.Item(132, 8).Value = 4360.978  # H132: 4690.976 -> 4360.978
$ws.Cells.Item(132, 9).Value = 4326.225  # I132: 4707.3613 -> 4326.225
$ws.Cells.Item(132, 11).Value = 12978.675  # K132: 14122.0839 -> 12978.675
$ws.Cells.Item(132, 13).Value = -10448.675  # M132: -11592.0839 -> -10448.675
# Row 139 (Leve Item ID 42306)
$ws.Cells.Item(139, 8).Value = 103826.11  # H139: 105476.43 -> 103826.11
$ws.Cells.Item(139, 10).Value = 103826.11  # J139: 105476.43 -> 103826.11
$ws.Cells.Item(139, 12).Value = 103826.11  # L139: 105476.43 -> 103826.11
$ws.Cells.Item(139, 14).Value = -114106.11  # N139: -115756.43 -> -114106.11

# ==== ARM sheet ====
$ws = $wb.Worksheets.Item("ARM")
# Row 23 (Leve Item ID 2236)
$ws.Cells.Item(23, 8).Value = 0  # H23: 2999 -> 0
$ws.Cells.Item(23, 10).Value = 0  # J23: 2999 -> 0
$ws.Cells.Item(23, 12).Value = 0  # L23: 2999 -> 0
$ws.Cells.Item(23, 14).ClearContents()  # N23: -3517 -> (cell removed)
# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 6351.879  # H32: 6437.484 -> 6351.879
$ws.Cells.Item(32, 9).Value = 6283.0312  # I32: 6437.484 -> 6283.0312
$ws.Cells.Item(32, 10).Value = 8555  # J32: 0 -> 8555
$ws.Cells.Item(32, 11).Value = 6283.0312  # K32: 6437.484 -> 6283.0312
$ws.Cells.Item(32, 12).Value = 8555  # L32: 0 -> 8555
$ws.Cells.Item(32, 13).Value = -5996.0312  # M32: -6150.484 -> -5996.0312
$ws.Cells.Item(32, 14).Value = -9129  # N32: None -> -9129
# Row 34 (Leve Item ID 2753)
$ws.Cells.Item(34, 8).Value = 184500  # H34: 223000 -> 184500
# Row 96 (Leve Item ID 18207)
$ws.Cells.Item(96, 8).Value = 0  # H96: 30000 -> 0
$ws.Cells.Item(96, 10).Value = 0  # J96: 30000 -> 0
$ws.Cells.Item(96, 12).Value = 0  # L96: 30000 -> 0
$ws.Cells.Item(96, 14).ClearContents()  # N96: -35492 -> (cell removed)
# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 981515.2  # H122: 1027157.44 -> 981515.2
$ws.Cells.Item(122, 9).Value = 2984.2258  # I122: 3077.0334 -> 2984.2258
$ws.Cells.Item(122, 10).Value = 3148262.2  # J122: 3390419.8 -> 3148262.2
$ws.Cells.Item(122, 11).Value = 8952.6774  # K122: 9231.100199999999 -> 8952.6774
$ws.Cells.Item(122, 12).Value = 9444786.600000001  # L122: 10171259.4 -> 9444786.600000001
$ws.Cells.Item(122, 13).Value = -6502.6774  # M122: -6781.100199999999 -> -6502.6774
$ws.Cells.Item(122, 14).Value = -9449686.600000001  # N122: -10176159.4 -> -9449686.600000001
# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 7746.394  # H132: 8078.484 -> 7746.394
$ws.Cells.Item(132, 9).Value = 9314.65  # I132: 10060.833 -> 9314.65
$ws.Cells.Item(132, 11).Value = 27943.95  # K132: 30182.499 -> 27943.95
$ws.Cells.Item(132, 13).Value = -25413.95  # M132: -27652.499 -> -25413.95

# ==== BSM sheet ====
$ws = $wb.Worksheets.Item("BSM")
# Row 25 (Leve Item ID 2370)
$ws.Cells.Item(25, 8).Value = 2302.3333  # H25: 6151 -> 2302.3333
$ws.Cells.Item(25, 9).Value = 456.5  # I25: 3637.3333 -> 456.5
$ws.Cells.Item(25, 10).Value = 5994  # J25: 8664.666999999999 -> 5994
$ws.Cells.Item(25, 11).Value = 456.5  # K25: 3637.3333 -> 456.5
$ws.Cells.Item(25, 12).Value = 5994  # L25: 8664.666999999999 -> 5994
$ws.Cells.Item(25, 13).Value = -221.5  # M25: -3402.3333 -> -221.5
$ws.Cells.Item(25, 14).Value = -6464  # N25: -9134.666999999999 -> -6464
# Row 86 (Leve Item ID 12526)
$ws.Cells.Item(86, 8).Value = 4268.7095  # H86: 4655 -> 4268.7095
$ws.Cells.Item(86, 9).Value = 4914.7915  # I86: 5293.864 -> 4914.7915
$ws.Cells.Item(86, 10).Value = 2053.5715  # J86: 2312.5 -> 2053.5715
$ws.Cells.Item(86, 11).Value = 4914.7915  # K86: 5293.864 -> 4914.7915
$ws.Cells.Item(86, 12).Value = 2053.5715  # L86: 2312.5 -> 2053.5715
$ws.Cells.Item(86, 13).Value = -3791.7915  # M86: -4170.864 -> -3791.7915
$ws.Cells.Item(86, 14).Value = -4299.5715  # N86: -4558.5 -> -4299.5715
# Row 89 (Leve Item ID 12526)
$ws.Cells.Item(89, 8).Value = 4268.7095  # H89: 4655 -> 4268.7095
$ws.Cells.Item(89, 9).Value = 4914.7915  # I89: 5293.864 -> 4914.7915
$ws.Cells.Item(89, 10).Value = 2053.5715  # J89: 2312.5 -> 2053.5715
$ws.Cells.Item(89, 11).Value = 24573.9575  # K89: 26469.32 -> 24573.9575
$ws.Cells.Item(89, 12).Value = 10267.8575  # L89: 11562.5 -> 10267.8575
$ws.Cells.Item(89, 13).Value = -18957.9575  # M89: -20853.32 -> -18957.9575
$ws.Cells.Item(89, 14).Value = -21499.8575  # N89: -22794.5 -> -21499.8575
# Row 94 (Leve Item ID 19939)
$ws.Cells.Item(94, 8).Value = 9099.257  # H94: 9636.727999999999 -> 9099.257
$ws.Cells.Item(94, 9).Value = 11529.76  # I94: 12512.261 -> 11529.76
$ws.Cells.Item(94, 11).Value = 11529.76  # K94: 12512.261 -> 11529.76
$ws.Cells.Item(94, 13).Value = -11078.76  # M94: -12061.261 -> -11078.76
# Row 107 (Leve Item ID 27706)
$ws.Cells.Item(107, 8).Value = 2693.3333  # H107: 3237.1428 -> 2693.3333
$ws.Cells.Item(107, 9).Value = 2530.625  # I107: 3110.8333 -> 2530.625
$ws.Cells.Item(107, 11).Value = 2530.625  # K107: 3110.8333 -> 2530.625
$ws.Cells.Item(107, 13).Value = -610.625  # M107: -1190.8333 -> -610.625

# ==== CRP sheet ====
$ws = $wb.Worksheets.Item("CRP")
# Row 88 (Leve Item ID 10608)
$ws.Cells.Item(88, 8).Value = 34370.223  # H88: 31646.6 -> 34370.223
$ws.Cells.Item(88, 9).Value = 24937.5  # I88: 33000 -> 24937.5
$ws.Cells.Item(88, 10).Value = 37065.285  # J88: 31308.25 -> 37065.285
$ws.Cells.Item(88, 11).Value = 24937.5  # K88: 33000 -> 24937.5
$ws.Cells.Item(88, 12).Value = 37065.285  # L88: 31308.25 -> 37065.285
$ws.Cells.Item(88, 13).Value = -24531.5  # M88: -32594 -> -24531.5
$ws.Cells.Item(88, 14).Value = -37877.285  # N88: -32120.25 -> -37877.285
# Row 91 (Leve Item ID 10608)
$ws.Cells.Item(91, 8).Value = 34370.223  # H91: 31646.6 -> 34370.223
$ws.Cells.Item(91, 9).Value = 24937.5  # I91: 33000 -> 24937.5
$ws.Cells.Item(91, 10).Value = 37065.285  # J91: 31308.25 -> 37065.285
$ws.Cells.Item(91, 11).Value = 24937.5  # K91: 33000 -> 24937.5
$ws.Cells.Item(91, 12).Value = 37065.285  # L91: 31308.25 -> 37065.285
$ws.Cells.Item(91, 13).Value = -23533.5  # M91: -31596 -> -23533.5
$ws.Cells.Item(91, 14).Value = -39873.285  # N91: -34116.25 -> -39873.285
# Row 99 (Leve Item ID 36198)
$ws.Cells.Item(99, 8).Value = 160964.47  # H99: 171590.47 -> 160964.47
$ws.Cells.Item(99, 9).Value = 360827.16  # I99: 420635.9 -> 360827.16
$ws.Cells.Item(99, 10).Value = 5515.722  # J99: 5560.1665 -> 5515.722
$ws.Cells.Item(99, 11).Value = 360827.16  # K99: 420635.9 -> 360827.16
$ws.Cells.Item(99, 12).Value = 5515.722  # L99: 5560.1665 -> 5515.722
$ws.Cells.Item(99, 13).Value = -359329.16  # M99: -419137.9 -> -359329.16
$ws.Cells.Item(99, 14).Value = -8511.722  # N99: -8556.166499999999 -> -8511.722
# Row 105 (Leve Item ID 19928)
$ws.Cells.Item(105, 8).Value = 113258.266  # H105: 119494.836 -> 113258.266
$ws.Cells.Item(105, 9).Value = 142380.53  # I105: 152479.14 -> 142380.53
$ws.Cells.Item(105, 11).Value = 142380.53  # K105: 152479.14 -> 142380.53
$ws.Cells.Item(105, 13).Value = -140633.53  # M105: -150732.14 -> -140633.53
# Row 107 (Leve Item ID 27689)
$ws.Cells.Item(107, 8).Value = 8932.223  # H107: 9272.346 -> 8932.223
$ws.Cells.Item(107, 9).Value = 9622.799999999999  # I107: 10020.042 -> 9622.799999999999
$ws.Cells.Item(107, 11).Value = 9622.799999999999  # K107: 10020.042 -> 9622.799999999999
$ws.Cells.Item(107, 13).Value = -7702.799999999999  # M107: -8100.041999999999 -> -7702.799999999999
# Row 121 (Leve Item ID 27227)
$ws.Cells.Item(121, 8).Value = 0  # H121: 74925 -> 0
$ws.Cells.Item(121, 10).Value = 0  # J121: 74925 -> 0
$ws.Cells.Item(121, 12).Value = 0  # L121: 74925 -> 0
$ws.Cells.Item(121, 14).ClearContents()  # N121: -77545 -> (cell removed)
# Row 122 (Leve Item ID 36196)
$ws.Cells.Item(122, 8).Value = 7364.6  # H122: 7418.8 -> 7364.6
$ws.Cells.Item(122, 9).Value = 11151.333  # I122: 11241.667 -> 11151.333
$ws.Cells.Item(122, 11).Value = 33453.999  # K122: 33725.001 -> 33453.999
$ws.Cells.Item(122, 13).Value = -31003.999  # M122: -31275.001 -> -31003.999
# Row 126 (Leve Item ID 36198)
$ws.Cells.Item(126, 8).Value = 160964.47  # H126: 171590.47 -> 160964.47
$ws.Cells.Item(126, 9).Value = 360827.16  # I126: 420635.9 -> 360827.16
$ws.Cells.Item(126, 10).Value = 5515.722  # J126: 5560.1665 -> 5515.722
$ws.Cells.Item(126, 11).Value = 1082481.48  # K126: 1261907.7 -> 1082481.48
$ws.Cells.Item(126, 12).Value = 16547.166  # L126: 16680.4995 -> 16547.166
$ws.Cells.Item(126, 13).Value = -1080011.48  # M126: -1259437.7 -> -1080011.48
$ws.Cells.Item(126, 14).Value = -21487.166  # N126: -21620.4995 -> -21487.166
# Row 134 (Leve Item ID 44020)
$ws.Cells.Item(134, 8).Value = 4677.4443  # H134: 3748.6667 -> 4677.4443
$ws.Cells.Item(134, 9).Value = 5683.3335  # I134: 4109.6665 -> 5683.3335
$ws.Cells.Item(134, 11).Value = 17050.0005  # K134: 12328.9995 -> 17050.0005
$ws.Cells.Item(134, 13).Value = -14515.0005  # M134: -9793.999500000002 -> -14515.0005

# ==== CUL sheet ====
$ws = $wb.Worksheets.Item("CUL")
# Row 70 (Leve Item ID 12867)
$ws.Cells.Item(70, 8).Value = 0  # H70: 2800 -> 0
$ws.Cells.Item(70, 9).Value = 0  # I70: 2800 -> 0
$ws.Cells.Item(70, 11).Value = 0  # K70: 8400 -> 0
$ws.Cells.Item(70, 13).ClearContents()  # M70: -8085 -> (cell removed)
# Row 73 (Leve Item ID 12867)
$ws.Cells.Item(73, 8).Value = 0  # H73: 2800 -> 0
$ws.Cells.Item(73, 9).Value = 0  # I73: 2800 -> 0
$ws.Cells.Item(73, 11).Value = 0  # K73: 8400 -> 0
$ws.Cells.Item(73, 13).ClearContents()  # M73: -7308 -> (cell removed)

# ==== GSM sheet ====
$ws = $wb.Worksheets.Item("GSM")
# Row 102 (Leve Item ID 36169)
$ws.Cells.Item(102, 8).Value = 7065.28  # H102: 4548.4146 -> 7065.28
$ws.Cells.Item(102, 9).Value = 10962.429  # I102: 5444.2334 -> 10962.429
$ws.Cells.Item(102, 11).Value = 10962.429  # K102: 5444.2334 -> 10962.429
$ws.Cells.Item(102, 13).Value = -9340.429  # M102: -3822.2334 -> -9340.429
# Row 113 (Leve Item ID 27710)
$ws.Cells.Item(113, 8).Value = 11868.066  # H113: 11939.4 -> 11868.066
$ws.Cells.Item(113, 9).Value = 13586.308  # I113: 13668.615 -> 13586.308
$ws.Cells.Item(113, 11).Value = 13586.308  # K113: 13668.615 -> 13586.308
$ws.Cells.Item(113, 13).Value = -11416.308  # M113: -11498.615 -> -11416.308
# Row 122 (Leve Item ID 36182)
$ws.Cells.Item(122, 8).Value = 5176.3125  # H122: 5267.1704 -> 5176.3125
$ws.Cells.Item(122, 10).Value = 13968.333  # J122: 15601.125 -> 13968.333
$ws.Cells.Item(122, 12).Value = 41904.999  # L122: 46803.375 -> 41904.999
$ws.Cells.Item(122, 14).Value = -46804.999  # N122: -51703.375 -> -46804.999
# Row 126 (Leve Item ID 36184)
$ws.Cells.Item(126, 8).Value = 15350.5625  # H126: 11564.228 -> 15350.5625
$ws.Cells.Item(126, 9).Value = 23234.834  # I126: 17666.875 -> 23234.834
$ws.Cells.Item(126, 10).Value = 10620  # J126: 8077 -> 10620
$ws.Cells.Item(126, 11).Value = 69704.50199999999  # K126: 53000.625 -> 69704.50199999999
$ws.Cells.Item(126, 12).Value = 31860  # L126: 24231 -> 31860
$ws.Cells.Item(126, 13).Value = -67234.50199999999  # M126: -50530.625 -> -67234.50199999999
$ws.Cells.Item(126, 14).Value = -36800  # N126: -29171 -> -36800
# Row 132 (Leve Item ID 44008)
$ws.Cells.Item(132, 8).Value = 1870.875  # H132: 1861 -> 1870.875
$ws.Cells.Item(132, 9).Value = 1762.2667  # I132: 1758.5625 -> 1762.2667
$ws.Cells.Item(132, 11).Value = 5286.800099999999  # K132: 5275.6875 -> 5286.800099999999
$ws.Cells.Item(132, 13).Value = -2756.800099999999  # M132: -2745.6875 -> -2756.800099999999

# ==== LTW sheet ====
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Cells.Item(22, 8).Value = 12353.556  # H22: 11814.211 -> 12353.556
$ws.Cells.Item(22, 10).Value = 1981.5  # J22: 1995.3334 -> 1981.5
$ws.Cells.Item(22, 12).Value = 1981.5  # L22: 1995.3334 -> 1981.5
$ws.Cells.Item(22, 14).Value = -2571.5  # N22: -2585.3334 -> -2571.5
# Row 27 (Leve Item ID 5277)
$ws.Cells.Item(27, 8).Value = 12353.556  # H27: 11814.211 -> 12353.556
$ws.Cells.Item(27, 10).Value = 1981.5  # J27: 1995.3334 -> 1981.5
$ws.Cells.Item(27, 12).Value = 1981.5  # L27: 1995.3334 -> 1981.5
$ws.Cells.Item(27, 14).Value = -2195.5  # N27: -2209.3334 -> -2195.5
# Row 29 (Leve Item ID 3576)
$ws.Cells.Item(29, 8).Value = 2999.75  # H29: 3666.6667 -> 2999.75
$ws.Cells.Item(29, 9).Value = 999.6667  # I29: 1000 -> 999.6667
$ws.Cells.Item(29, 11).Value = 999.6667  # K29: 1000 -> 999.6667
$ws.Cells.Item(29, 13).Value = -704.6667  # M29: -705 -> -704.6667
# Row 68 (Leve Item ID 12563)
$ws.Cells.Item(68, 8).Value = 5257.8  # H68: 6239.7144 -> 5257.8
$ws.Cells.Item(68, 9).Value = 2296.3333  # I68: 2219.5 -> 2296.3333
$ws.Cells.Item(68, 10).Value = 9700  # J68: 11600 -> 9700
$ws.Cells.Item(68, 11).Value = 2296.3333  # K68: 2219.5 -> 2296.3333
$ws.Cells.Item(68, 12).Value = 9700  # L68: 11600 -> 9700
$ws.Cells.Item(68, 13).Value = -1547.3333  # M68: -1470.5 -> -1547.3333
$ws.Cells.Item(68, 14).Value = -11198  # N68: -13098 -> -11198
# Row 71 (Leve Item ID 12563)
$ws.Cells.Item(71, 8).Value = 5257.8  # H71: 6239.7144 -> 5257.8
$ws.Cells.Item(71, 9).Value = 2296.3333  # I71: 2219.5 -> 2296.3333
$ws.Cells.Item(71, 10).Value = 9700  # J71: 11600 -> 9700
$ws.Cells.Item(71, 11).Value = 11481.6665  # K71: 11097.5 -> 11481.6665
$ws.Cells.Item(71, 12).Value = 48500  # L71: 58000 -> 48500
$ws.Cells.Item(71, 13).Value = -7737.666499999999  # M71: -7353.5 -> -7737.666499999999
$ws.Cells.Item(71, 14).Value = -55988  # N71: -65488 -> -55988
# Row 122 (Leve Item ID 36247)
$ws.Cells.Item(122, 8).Value = 4071.8647  # H122: 4111.081 -> 4071.8647
$ws.Cells.Item(122, 9).Value = 5474.5454  # I122: 5264.25 -> 5474.5454
$ws.Cells.Item(122, 10).Value = 3478.423  # J122: 3557.56 -> 3478.423
$ws.Cells.Item(122, 11).Value = 16423.6362  # K122: 15792.75 -> 16423.6362
$ws.Cells.Item(122, 12).Value = 10435.269  # L122: 10672.68 -> 10435.269
$ws.Cells.Item(122, 13).Value = -13973.6362  # M122: -13342.75 -> -13973.6362
$ws.Cells.Item(122, 14).Value = -15335.269  # N122: -15572.68 -> -15335.269
# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 750240.9  # H132: 833357.75 -> 750240.9
$ws.Cells.Item(132, 9).Value = 1493983.1  # I132: 1866931.5 -> 1493983.1
$ws.Cells.Item(132, 11).Value = 4481949.300000001  # K132: 5600794.5 -> 4481949.300000001
$ws.Cells.Item(132, 13).Value = -4479419.300000001  # M132: -5598264.5 -> -4479419.300000001

# ==== WVR sheet ====
$ws = $wb.Worksheets.Item("WVR")
# Row 122 (Leve Item ID 36208)
$ws.Cells.Item(122, 8).Value = 3769.2114  # H122: 3835.2354 -> 3769.2114
$ws.Cells.Item(122, 9).Value = 1584.6154  # I122: 1615.7368 -> 1584.6154
$ws.Cells.Item(122, 11).Value = 4753.8462  # K122: 4847.2104 -> 4753.8462
$ws.Cells.Item(122, 13).Value = -2303.8462  # M122: -2397.2104 -> -2303.8462
# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, 8).Value = 6781.644  # H132: 6605.92 -> 6781.644
$ws.Cells.Item(132, 9).Value = 8067.423  # I132: 7775.7407 -> 8067.423
$ws.Cells.Item(132, 11).Value = 24202.269  # K132: 23327.2221 -> 24202.269
$ws.Cells.Item(132, 13).Value = -21672.269  # M132: -20797.2221 -> -21672.269
# Row 136 (Leve Item ID 44031)
$ws.Cells.Item(136, 8).Value = 246253.23  # H136: 254313.53 -> 246253.23
$ws.Cells.Item(136, 9).Value = 343516.66  # I136: 351324.3 -> 343516.66
$ws.Cells.Item(136, 10).Value = 3094.6667  # J136: 3226.7646 -> 3094.6667
$ws.Cells.Item(136, 11).Value = 1030549.98  # K136: 1053972.9 -> 1030549.98
$ws.Cells.Item(136, 12).Value = 9284.000100000001  # L136: 9680.293799999999 -> 9284.000100000001
$ws.Cells.Item(136, 13).Value = -1027999.98  # M136: -1051422.9 -> -1027999.98
$ws.Cells.Item(136, 14).Value = -14384.0001  # N136: -14780.2938 -> -14384.0001
